$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 29 (2025Q3) metrics per the diff
$ws.Range("C29").Value = 176
$ws.Range("D29").Value = 25
$ws.Range("E29").Value = 151
$ws.Range("F29").Value = 4.302925989672977
